# Auto-generated Excel COM-interop script to apply the profit-table refresh
# produced by the scheduled Sheets runner. For each affected sheet, update the
# changed numeric cells (currentAveragePrice*, LevePrice*, LeveProfit*) to match
# the latest market-board snapshot; remove LeveProfitHQ cells that no longer apply.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 44750
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 44750
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 134250
$ws.Range("N57").Value = -135248
$ws.Range("H64").Value = 3778.2222
$ws.Range("I64").Value = 3667.3333
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3667.3333
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3419.3333
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 3778.2222
$ws.Range("I67").Value = 3667.3333
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3667.3333
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2809.3333
$ws.Range("N67").Value = -5716
$ws.Range("H74").Value = 3584.4736
$ws.Range("I74").Value = 3375.25
$ws.Range("J74").Value = 3943.1428
$ws.Range("K74").Value = 3375.25
$ws.Range("L74").Value = 3943.1428
$ws.Range("M74").Value = -2439.25
$ws.Range("N74").Value = -5815.1428
$ws.Range("H77").Value = 3584.4736
$ws.Range("I77").Value = 3375.25
$ws.Range("J77").Value = 3943.1428
$ws.Range("K77").Value = 16876.25
$ws.Range("L77").Value = 19715.714
$ws.Range("M77").Value = -12196.25
$ws.Range("N77").Value = -29075.714
$ws.Range("H132").Value = 2008.1111
$ws.Range("I132").Value = 1400.5416
$ws.Range("J132").Value = 6868.6665
$ws.Range("K132").Value = 4201.6248
$ws.Range("L132").Value = 20605.9995
$ws.Range("M132").Value = -1671.6248
$ws.Range("N132").Value = -25665.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3862.5
$ws.Range("I63").Value = 2725
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2725
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -2039
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 3862.5
$ws.Range("I66").Value = 2725
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 13625
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -10193
$ws.Range("N66").Value = -31864
$ws.Range("H102").Value = 2705
$ws.Range("I102").Value = 2410
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2410
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -788
$ws.Range("N102").Value = -6244
$ws.Range("H133").Value = 27997.75
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 27997.75
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 27997.75
$ws.Range("N133").Value = -33057.75
$ws.Range("H134").Value = 30666.334
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 30666.334
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 30666.334
$ws.Range("N134").Value = -40806.334
$ws.Range("H135").Value = 28760
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 28760
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 28760
$ws.Range("N135").Value = -38900
$ws.Range("H137").Value = 54711.285
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 54711.285
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 54711.285
$ws.Range("N137").Value = -64911.285
$ws.Range("H138").Value = 35000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 35000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 35000
$ws.Range("N138").Value = -45280
$ws.Range("H139").Value = 37000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 37000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 37000
$ws.Range("N139").Value = -47280
$ws.Range("H140").Value = 34299.75
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 34299.75
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 34299.75
$ws.Range("N140").Value = -44659.75
$ws.Range("H141").Value = 45150
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 45150
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 45150
$ws.Range("N141").Value = -55510

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2136.5264
$ws.Range("I86").Value = 2049.9167
$ws.Range("J86").Value = 2285
$ws.Range("K86").Value = 2049.9167
$ws.Range("L86").Value = 2285
$ws.Range("M86").Value = -926.9167000000002
$ws.Range("N86").Value = -4531
$ws.Range("H89").Value = 2136.5264
$ws.Range("I89").Value = 2049.9167
$ws.Range("J89").Value = 2285
$ws.Range("K89").Value = 10249.5835
$ws.Range("L89").Value = 11425
$ws.Range("M89").Value = -4633.583500000001
$ws.Range("N89").Value = -22657
$ws.Range("H94").Value = 822.7059
$ws.Range("I94").Value = 778
$ws.Range("J94").Value = 930
$ws.Range("K94").Value = 778
$ws.Range("L94").Value = 930
$ws.Range("M94").Value = -327
$ws.Range("N94").Value = -1832
$ws.Range("H103").Value = 27771.334
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 27771.334
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 27771.334
$ws.Range("N103").Value = -30115.334
$ws.Range("H134").Value = 4906.204
$ws.Range("I134").Value = 5031.2646
$ws.Range("J134").Value = 4622.7334
$ws.Range("K134").Value = 15093.7938
$ws.Range("L134").Value = 13868.2002
$ws.Range("M134").Value = -12558.7938
$ws.Range("N134").Value = -18938.2002
$ws.Range("H139").Value = 51445
$ws.Range("I139").Value = 40000
$ws.Range("J139").Value = 55260
$ws.Range("K139").Value = 40000
$ws.Range("L139").Value = 55260
$ws.Range("M139").Value = -34860
$ws.Range("N139").Value = -65540
$ws.Range("H140").Value = 39750
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 39750
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 39750
$ws.Range("N140").Value = -50110

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4836.8667
$ws.Range("I62").Value = 2472.5
$ws.Range("J62").Value = 7539
$ws.Range("K62").Value = 2472.5
$ws.Range("L62").Value = 7539
$ws.Range("M62").Value = -1848.5
$ws.Range("N62").Value = -8787
$ws.Range("H65").Value = 4836.8667
$ws.Range("I65").Value = 2472.5
$ws.Range("J65").Value = 7539
$ws.Range("K65").Value = 12362.5
$ws.Range("L65").Value = 37695
$ws.Range("M65").Value = -9242.5
$ws.Range("N65").Value = -43935
$ws.Range("H132").Value = 2573.2593
$ws.Range("I132").Value = 2292.2727
$ws.Range("J132").Value = 3809.6
$ws.Range("K132").Value = 6876.8181
$ws.Range("L132").Value = 11428.8
$ws.Range("M132").Value = -4346.8181
$ws.Range("N132").Value = -16488.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4333.7144
$ws.Range("I64").Value = 4087.2
$ws.Range("J64").Value = 4950
$ws.Range("K64").Value = 12261.6
$ws.Range("L64").Value = 14850
$ws.Range("M64").Value = -11991.6
$ws.Range("N64").Value = -15390
$ws.Range("H67").Value = 4333.7144
$ws.Range("I67").Value = 4087.2
$ws.Range("J67").Value = 4950
$ws.Range("K67").Value = 12261.6
$ws.Range("L67").Value = 14850
$ws.Range("M67").Value = -11325.6
$ws.Range("N67").Value = -16722
$ws.Range("H70").Value = 4951.952
$ws.Range("I70").Value = 3071.1
$ws.Range("J70").Value = 6661.8184
$ws.Range("K70").Value = 9213.299999999999
$ws.Range("L70").Value = 19985.4552
$ws.Range("M70").Value = -8898.299999999999
$ws.Range("N70").Value = -20615.4552
$ws.Range("H73").Value = 4951.952
$ws.Range("I73").Value = 3071.1
$ws.Range("J73").Value = 6661.8184
$ws.Range("K73").Value = 9213.299999999999
$ws.Range("L73").Value = 19985.4552
$ws.Range("M73").Value = -8121.299999999999
$ws.Range("N73").Value = -22169.4552
$ws.Range("H76").Value = 3450
$ws.Range("I76").Value = 1000
$ws.Range("J76").Value = 5900
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 17700
$ws.Range("M76").Value = -2617
$ws.Range("N76").Value = -18466
$ws.Range("H79").Value = 3450
$ws.Range("I79").Value = 1000
$ws.Range("J79").Value = 5900
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 17700
$ws.Range("M79").Value = -1674
$ws.Range("N79").Value = -20352
$ws.Range("H131").Value = 914.94666
$ws.Range("I131").Value = 473.375
$ws.Range("J131").Value = 967.6716300000001
$ws.Range("K131").Value = 1420.125
$ws.Range("L131").Value = 2903.01489
$ws.Range("M131").Value = 3619.875
$ws.Range("N131").Value = -12983.01489

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3770.25
$ws.Range("I80").Value = 2659.1177
$ws.Range("J80").Value = 10066.667
$ws.Range("K80").Value = 2659.1177
$ws.Range("L80").Value = 10066.667
$ws.Range("M80").Value = -1661.1177
$ws.Range("N80").Value = -12062.667
$ws.Range("H83").Value = 3770.25
$ws.Range("I83").Value = 2659.1177
$ws.Range("J83").Value = 10066.667
$ws.Range("K83").Value = 13295.5885
$ws.Range("L83").Value = 50333.335
$ws.Range("M83").Value = -8303.588499999998
$ws.Range("N83").Value = -60317.335
$ws.Range("H122").Value = 1225.8667
$ws.Range("I122").Value = 1320.8889
$ws.Range("J122").Value = 1083.3334
$ws.Range("K122").Value = 3962.6667
$ws.Range("L122").Value = 3250.0002
$ws.Range("M122").Value = -1512.6667
$ws.Range("H132").Value = 4152.9414
$ws.Range("I132").Value = 4290.2
$ws.Range("J132").Value = 3956.8572
$ws.Range("K132").Value = 12870.6
$ws.Range("L132").Value = 11870.5716
$ws.Range("M132").Value = -10340.6
$ws.Range("N132").Value = -16930.5716
$ws.Range("H138").Value = 34300
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 34300
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 34300
$ws.Range("N138").Value = -44580
$ws.Range("H140").Value = 39739.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 39739.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 39739.5
$ws.Range("N140").Value = -50099.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1159.619
$ws.Range("I82").Value = 1038.4
$ws.Range("J82").Value = 1197.5
$ws.Range("K82").Value = 1038.4
$ws.Range("L82").Value = 1197.5
$ws.Range("M82").Value = -677.4000000000001
$ws.Range("N82").Value = -1919.5
$ws.Range("H85").Value = 1159.619
$ws.Range("I85").Value = 1038.4
$ws.Range("J85").Value = 1197.5
$ws.Range("K85").Value = 1038.4
$ws.Range("L85").Value = 1197.5
$ws.Range("M85").Value = 209.5999999999999
$ws.Range("N85").Value = -3693.5
$ws.Range("H132").Value = 2383.3713
$ws.Range("I132").Value = 2402.7112
$ws.Range("J132").Value = 2348.56
$ws.Range("K132").Value = 7208.133600000001
$ws.Range("L132").Value = 7045.68
$ws.Range("M132").Value = -4678.133600000001
$ws.Range("N132").Value = -12105.68

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("H132").Value = 4483.5557
$ws.Range("I132").Value = 5183.5
$ws.Range("J132").Value = 3083.6667
$ws.Range("K132").Value = 15550.5
$ws.Range("L132").Value = 9251.000100000001
$ws.Range("M132").Value = -13020.5
$ws.Range("N132").Value = -14311.0001
$ws.Range("H136").Value = 4484.0557
$ws.Range("I136").Value = 8008.2856
$ws.Range("J136").Value = 2241.3635
$ws.Range("K136").Value = 24024.8568
$ws.Range("L136").Value = 6724.0905
$ws.Range("M136").Value = -21474.8568
$ws.Range("N136").Value = -11824.0905
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()
